# Team 44 Info.xlsx — add Tyler Dalke's row to the team roster table.
#
# Layout is A:Name, B:Preferred Email (hyperlinked "Hyperlink" style like
# the rows above it), C:Country, D:Time Zone. New data goes in row 4,
# immediately below the last existing person (row 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Tyler Dalke"
$ws.Range("B4").Value = "TylerADalke22@hotmail.com"
$ws.Range("C4").Value = "Canada"
$ws.Range("D4").Value = "UTC-5"

# Match the look of the other "Preferred Email" cells (B2, B3): an actual
# mailto hyperlink plus the underlined "Hyperlink" cell style. (Style is
# applied after Hyperlinks.Add so it reuses the same style slot as B2/B3
# instead of allocating a fresh one.)
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:TylerADalke22@hotmail.com")
$ws.Range("B4").Style = "Hyperlink"

# Return the selection to the top-left cell (matches the saved file no
# longer pointing the cursor at the old, now-stale E3 selection).
[void]$ws.Range("A1").Select()
